$d = $word.ActiveDocument

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# --- Step 1: expand the single paragraph into 4 paragraphs ---
# (paragraph 1 keeps "hello" and the _GoBack bookmark for now; the new
#  paragraphs 2-4 start out empty)
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()
$p1.Range.InsertParagraphAfter()
$p1.Range.InsertParagraphAfter()

# --- Step 2: drop the _GoBack bookmark; it gets re-created at the very
#     end of the document once the final text is in place ---
$d.Bookmarks("_GoBack").Delete()

# --- Step 3: paragraph 1 - split "hello" into "H" / "ello" / ". How are
#     you doing. " as three separate runs ---
$p1 = $d.Paragraphs(1)
$r1 = $d.Range($p1.Range.Start, $p1.Range.End)
$xml1 = "<w:p xmlns:w='$wNs'>" +
        "<w:r><w:t>H</w:t></w:r>" +
        "<w:r><w:t>ello</w:t></w:r>" +
        "<w:r><w:t xml:space='preserve'>. How are you doing. </w:t></w:r>" +
        "</w:p>"
$r1.InsertXML($xml1)

# --- Step 4: paragraph 2 - "Are you " + [gramStart] "Ok." [gramEnd] + " " ---
$p2 = $d.Paragraphs(2)
$r2 = $d.Range($p2.Range.Start, $p2.Range.End)
$xml2 = "<w:p xmlns:w='$wNs'>" +
        "<w:r><w:t xml:space='preserve'>Are you </w:t></w:r>" +
        "<w:proofErr w:type='gramStart'/>" +
        "<w:r><w:t>Ok.</w:t></w:r>" +
        "<w:proofErr w:type='gramEnd'/>" +
        "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
        "</w:p>"
$r2.InsertXML($xml2)

# --- Step 5: paragraph 3 stays a bare empty paragraph ---
$p3 = $d.Paragraphs(3)
$r3 = $d.Range($p3.Range.Start, $p3.Range.End)
$r3.InsertXML("<w:p xmlns:w='$wNs'/>")

# --- Step 6: paragraph 4 - final sentence ---
$p4 = $d.Paragraphs(4)
$r4 = $d.Range($p4.Range.Start, $p4.Range.End)
$xml4 = "<w:p xmlns:w='$wNs'>" +
        "<w:r><w:t>Let me know If there is any help I can do .</w:t></w:r>" +
        "</w:p>"
$r4.InsertXML($xml4)

# --- Step 7: re-add the _GoBack bookmark at the very end of paragraph 4 ---
$p4 = $d.Paragraphs(4)
$endPos = $p4.Range.End - 1
$bkRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bkRange)
